# September_Schedule.xlsx
# Fill in the "Flying/Not Flying" master table for the two remaining weeks
# (rows 6 and 7: weeks beginning 2022-09-19 and 2022-09-26).
#
# Row 6 alternates between the two activities across the Morning/Afternoon
# column pairs (B..K); row 7 is "Data analysis, report writing" for every
# Morning/Afternoon slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$concentrations = "Concentrations, preparing boxes to run"
$dataAnalysis   = "Data analysis, report writing"

# Row 6: alternating Concentrations / Data analysis across B6:K6
for ($col = 2; $col -le 11; $col++) {
    if (($col % 2) -eq 0) {
        $ws.Cells.Item(6, $col).Value = $concentrations
    } else {
        $ws.Cells.Item(6, $col).Value = $dataAnalysis
    }
}

# Row 7: Data analysis, report writing across the whole row B7:K7
for ($col = 2; $col -le 11; $col++) {
    $ws.Cells.Item(7, $col).Value = $dataAnalysis
}

# The rows now wrap onto several lines of text - resize them the way Excel
# would after the content grew.
$ws.Rows.Item(6).RowHeight = 59.7
$ws.Rows.Item(7).RowHeight = 48.05

# Leave the selection where the author's editing session ended up.
$ws.Range("F8").Select()
